# Apply the edit described by commit "LoopWithAndSplit and AndSplitWithLoop #422"
# Target sheet: "AndSplitWithLoop" (xl/worksheets/sheet5.xml)
#
# Before (rows 3-10):
#   3  AndSplit
#   4  Loop
#   5  Elementary | TestItem_AndSplit:0 | Loop
#   6  LoopEnd
#   7  Block
#   8  Elementary | TestItem_AndSplit:0 | Right
#   9  End
#   10 End
#
# After (rows 3-12): a new Block/End pair now wraps the Loop, and the
# routing expression of the Elementary step inside the loop becomes the
# new shared string "LeftInLoop" instead of "Loop".
#   3  AndSplit
#   4  Block
#   5  Loop
#   6  Elementary | TestItem_AndSplit:0 | LeftInLoop
#   7  LoopEnd
#   8  End
#   9  Block
#   10 Elementary | TestItem_AndSplit:0 | Right
#   11 End
#   12 End

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AndSplitWithLoop")

# Give the whole (now larger) A4:C12 block the same formatting as the
# existing data rows (row 3) -- this reuses the existing style instead of
# Excel's row-Insert() machinery, which tends to register spurious unused
# styles.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C12").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Clear out any stale values left over from the old, shorter table so every
# cell below ends up holding exactly the value we want (and nothing else).
$ws.Range("A4:C12").ClearContents()

$ws.Cells.Item(4, 1).Value2 = "Block"

$ws.Cells.Item(5, 1).Value2 = "Loop"

$ws.Cells.Item(6, 1).Value2 = "Elementary"
$ws.Cells.Item(6, 2).Value2 = "TestItem_AndSplit:0"
$ws.Cells.Item(6, 3).Value2 = "LeftInLoop"

$ws.Cells.Item(7, 1).Value2 = "LoopEnd"

$ws.Cells.Item(8, 1).Value2 = "End"

$ws.Cells.Item(9, 1).Value2 = "Block"

$ws.Cells.Item(10, 1).Value2 = "Elementary"
$ws.Cells.Item(10, 2).Value2 = "TestItem_AndSplit:0"
$ws.Cells.Item(10, 3).Value2 = "Right"

$ws.Cells.Item(11, 1).Value2 = "End"

$ws.Cells.Item(12, 1).Value2 = "End"

# Column C grew a touch wider to fit "LeftInLoop".
$ws.Columns.Item(3).ColumnWidth = 9.6

# The sheet is now shown with A8 selected instead of the old I32.
$ws.Range("A8").Select()
